{"js": "// Helper: search for `searchText` in the document body and replace every\n// match with `replacement`. Returns the number of matches processed.\nasync function replaceAll(searchText, replacement, options) {\n  const results = context.document.body.search(\n    searchText,\n    options || { matchCase: true }\n  );\n  results.load(\"text\");\n  await context.sync();\n\n  const count = results.items.length;\n  for (let i = 0; i < count; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return count;\n}\n\n// 1) Item description: \"MOVIL\" -> \"REFRIGERADOR LG\" (two occurrences; one run\n//    keeps a trailing space, the other doesn't - searching the bare word and\n//    replacing with the bare word leaves the surrounding whitespace runs\n//    untouched).\nawait replaceAll(\"MOVIL\", \"REFRIGERADOR LG\");\n\n// 2) The long composed date \"26 DE AGOSTO DE 2024\" -> \"15 DE ENERO DE 2025\"\n//    (single run). Do this before touching the standalone \"26\"/\"8\" runs so\n//    the generic whole-word searches below don't also see this one.\nawait replaceAll(\"26 DE AGOSTO DE 2024\", \"15 DE ENERO DE 2025\");\n\n// 3) Money amounts.\nawait replaceAll(\"$598.19\", \"$499.63\");\nawait replaceAll(\"$120.00\", \"$0\");\nawait replaceAll(\"$59.77\", \"$166.54\");\n\n// 4) Standalone \"8\" runs (plazo en meses / cuotas / garantia) -> \"3\".\n//    matchWholeWord avoids partial hits inside ids like \"98875666-3\".\nawait replaceAll(\"8\", \"3\", { matchCase: true, matchWholeWord: true });\n\n// 5) Standalone \"26\" runs (dia de pago) -> \"15\". The \"26\" that was part of\n//    \"26 DE AGOSTO DE 2024\" is already gone (step 2), so only the three\n//    standalone day-of-month runs remain.\nawait replaceAll(\"26\", \"15\", { matchCase: true, matchWholeWord: true });\n\n// 6) Notarial date pieces (day already handled above as part of the\n//    whole-word \"26\"; month + year remain).\nawait replaceAll(\"AGOSTO\", \"ENERO\");\nawait replaceAll(\" 2024\", \" 2025\");\n\n// 7) Tenant (\"arrendatario\") name and ID, repeated twice each in the\n//    notarial appearance clause.\nawait replaceAll(\"JAIME EDGARDO PALACIOS GARCIA\", \"JOSE ANTONIO PEREZ\");\nawait replaceAll(\"05682717-5\", \"87877868-7\");\n\n// 8) Guarantor (\"fiador\") name/address/ID fields are cleared out (emptied)\n//    - also repeated twice each.\nawait replaceAll(\"JORGE ANTONIO LOPEZ\", \"\");\nawait replaceAll(\"QUEZALTEPEQUE, LA LIBERTAD\", \"\");\nawait replaceAll(\"09388877-7\", \"\");\n", "ps1": "$d = $word.ActiveDocument\n\nfunction ReplaceAll {\n    param(\n        [string]$SearchText,\n        [string]$ReplaceText,\n        [bool]$WholeWord = $false\n    )\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $WholeWord\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $ReplaceText\n\n    # wdFindContinue = 1, Replace:=wdReplaceAll (2) so every match in the\n    # document gets replaced in one call.\n    $find.Execute($find.Text, $true, $WholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# 1) Item description: \"MOVIL\" -> \"REFRIGERADOR LG\" (two occurrences; one run\n#    keeps a trailing space, the other doesn't - searching the bare word and\n#    replacing with the bare word leaves the surrounding whitespace runs\n#    untouched).\nReplaceAll \"MOVIL\" \"REFRIGERADOR LG\" $false\n\n# 2) The long composed date \"26 DE AGOSTO DE 2024\" -> \"15 DE ENERO DE 2025\"\n#    (single run). Do this before touching the standalone \"26\"/\"8\" runs so\n#    the generic whole-word searches below don't also see this one.\nReplaceAll \"26 DE AGOSTO DE 2024\" \"15 DE ENERO DE 2025\" $false\n\n# 3) Money amounts.\nReplaceAll \"`$598.19\" \"`$499.63\" $false\nReplaceAll \"`$120.00\" \"`$0\" $false\nReplaceAll \"`$59.77\" \"`$166.54\" $false\n\n# 4) Standalone \"8\" runs (plazo en meses / cuotas / garantia) -> \"3\".\n#    MatchWholeWord avoids partial hits inside ids like \"98875666-3\".\nReplaceAll \"8\" \"3\" $true\n\n# 5) Standalone \"26\" runs (dia de pago) -> \"15\". The \"26\" that was part of\n#    \"26 DE AGOSTO DE 2024\" is already gone (step 2), so only the three\n#    standalone day-of-month runs remain.\nReplaceAll \"26\" \"15\" $true\n\n# 6) Notarial date pieces (day already handled above as part of the\n#    whole-word \"26\"; month + year remain).\nReplaceAll \"AGOSTO\" \"ENERO\" $false\nReplaceAll \" 2024\" \" 2025\" $false\n\n# 7) Tenant (\"arrendatario\") name and ID, repeated twice each in the\n#    notarial appearance clause.\nReplaceAll \"JAIME EDGARDO PALACIOS GARCIA\" \"JOSE ANTONIO PEREZ\" $false\nReplaceAll \"05682717-5\" \"87877868-7\" $false\n\n# 8) Guarantor (\"fiador\") name/address/ID fields are cleared out (emptied)\n#    - also repeated twice each.\nReplaceAll \"JORGE ANTONIO LOPEZ\" \"\" $false\nReplaceAll \"QUEZALTEPEQUE, LA LIBERTAD\" \"\" $false\nReplaceAll \"09388877-7\" \"\" $false\n"}
